# Update column F (dSF) values on Sheet1 to reflect the re-pulled/recalculated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = 6
$ws.Range("F8").Value = -4
$ws.Range("F10").Value = -4
$ws.Range("F12").Value = 9
$ws.Range("F13").Value = 4
$ws.Range("F15").Value = -1
